$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 1756.0714  # H111 was 1783.0358
$ws.Cells.Item(111, 10).Value = 1813.1154  # J111 was 1842.1538
$ws.Cells.Item(111, 12).Value = 5439.3462  # L111 was 5526.4614
$ws.Cells.Item(111, 14).Value = -11573.3462  # N111 was -11660.4614
$ws.Cells.Item(132, 8).Value = 9724.708000000001  # H132 was 10535.182
$ws.Cells.Item(132, 9).Value = 10108.131  # I132 was 10993.714
$ws.Cells.Item(132, 11).Value = 30324.393  # K132 was 32981.142
$ws.Cells.Item(132, 13).Value = -27794.393  # M132 was -30451.142
$ws.Cells.Item(141, 8).Value = 7911.5  # H141 was 7912.1665
$ws.Cells.Item(141, 9).Value = 7496.3335  # I141 was 7495.75
$ws.Cells.Item(141, 10).Value = 8326.666999999999  # J141 was 8745
$ws.Cells.Item(141, 11).Value = 22489.0005  # K141 was 22487.25
$ws.Cells.Item(141, 12).Value = 24980.001  # L141 was 26235
$ws.Cells.Item(141, 13).Value = -17309.0005  # M141 was -17307.25
$ws.Cells.Item(141, 14).Value = -35340.001  # N141 was -36595

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 0  # H122 was 1199.5
$ws.Cells.Item(122, 9).Value = 0  # I122 was 1199.5
$ws.Cells.Item(122, 11).Value = 0  # K122 was 3598.5
$ws.Cells.Item(122, 13).ClearContents()  # M122 was -1148.5
$ws.Cells.Item(138, 8).Value = 0  # H138 was 90429
$ws.Cells.Item(138, 10).Value = 0  # J138 was 90429
$ws.Cells.Item(138, 12).Value = 0  # L138 was 90429
$ws.Cells.Item(138, 14).ClearContents()  # N138 was -100709
$ws.Cells.Item(139, 8).Value = 0  # H139 was 25000
$ws.Cells.Item(139, 10).Value = 0  # J139 was 25000
$ws.Cells.Item(139, 12).Value = 0  # L139 was 25000
$ws.Cells.Item(139, 14).ClearContents()  # N139 was -35280

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(38, 8).Value = 750  # H38 was 2509.5
$ws.Cells.Item(38, 9).Value = 750  # I38 was 519
$ws.Cells.Item(38, 10).Value = 0  # J38 was 4500
$ws.Cells.Item(38, 11).Value = 750  # K38 was 519
$ws.Cells.Item(38, 12).Value = 0  # L38 was 4500
$ws.Cells.Item(38, 13).Value = -373  # M38 was -142
$ws.Cells.Item(38, 14).ClearContents()  # N38 was -5254
$ws.Cells.Item(46, 8).Value = 750  # H46 was 2509.5
$ws.Cells.Item(46, 9).Value = 750  # I46 was 519
$ws.Cells.Item(46, 10).Value = 0  # J46 was 4500
$ws.Cells.Item(46, 11).Value = 750  # K46 was 519
$ws.Cells.Item(46, 12).Value = 0  # L46 was 4500
$ws.Cells.Item(46, 13).Value = -539  # M46 was -308
$ws.Cells.Item(46, 14).ClearContents()  # N46 was -4922
$ws.Cells.Item(58, 8).Value = 1480.8182  # H58 was 2898.5454
$ws.Cells.Item(58, 9).Value = 1480.8182  # I58 was 1543.7778
$ws.Cells.Item(58, 10).Value = 0  # J58 was 8995
$ws.Cells.Item(58, 11).Value = 1480.8182  # K58 was 1543.7778
$ws.Cells.Item(58, 12).Value = 0  # L58 was 8995
$ws.Cells.Item(58, 13).Value = -1277.8182  # M58 was -1340.7778
$ws.Cells.Item(58, 14).ClearContents()  # N58 was -9401
$ws.Cells.Item(59, 8).Value = 42691.25  # H59 was 42441.25
$ws.Cells.Item(59, 10).Value = 42691.25  # J59 was 42441.25
$ws.Cells.Item(59, 12).Value = 42691.25  # L59 was 42441.25
$ws.Cells.Item(59, 14).Value = -44981.25  # N59 was -44731.25
$ws.Cells.Item(62, 8).Value = 1100  # H62 was 1250
$ws.Cells.Item(62, 9).Value = 1100  # I62 was 1250
$ws.Cells.Item(62, 11).Value = 1100  # K62 was 1250
$ws.Cells.Item(62, 13).Value = -476  # M62 was -626
$ws.Cells.Item(65, 8).Value = 1100  # H65 was 1250
$ws.Cells.Item(65, 9).Value = 1100  # I65 was 1250
$ws.Cells.Item(65, 11).Value = 5500  # K65 was 6250
$ws.Cells.Item(65, 13).Value = -2380  # M65 was -3130
$ws.Cells.Item(80, 8).Value = 47749.5  # H80 was 49999.5
$ws.Cells.Item(80, 10).Value = 47749.5  # J80 was 49999.5
$ws.Cells.Item(80, 12).Value = 47749.5  # L80 was 49999.5
$ws.Cells.Item(80, 14).Value = -49995.5  # N80 was -52245.5
$ws.Cells.Item(83, 8).Value = 47749.5  # H83 was 49999.5
$ws.Cells.Item(83, 10).Value = 47749.5  # J83 was 49999.5
$ws.Cells.Item(83, 12).Value = 143248.5  # L83 was 149998.5
$ws.Cells.Item(83, 14).Value = -154480.5  # N83 was -161230.5
$ws.Cells.Item(86, 8).Value = 4101.6665  # H86 was 4252.5
$ws.Cells.Item(86, 9).Value = 4101.6665  # I86 was 4252.5
$ws.Cells.Item(86, 11).Value = 4101.6665  # K86 was 4252.5
$ws.Cells.Item(86, 13).Value = -2978.6665  # M86 was -3129.5
$ws.Cells.Item(89, 8).Value = 4101.6665  # H89 was 4252.5
$ws.Cells.Item(89, 9).Value = 4101.6665  # I89 was 4252.5
$ws.Cells.Item(89, 11).Value = 20508.3325  # K89 was 21262.5
$ws.Cells.Item(89, 13).Value = -14892.3325  # M89 was -15646.5
$ws.Cells.Item(132, 8).Value = 1739.3334  # H132 was 2144.5454
$ws.Cells.Item(132, 9).Value = 1739.3334  # I132 was 2144.5454
$ws.Cells.Item(132, 11).Value = 5218.0002  # K132 was 6433.6362
$ws.Cells.Item(132, 13).Value = -2688.0002  # M132 was -3903.6362
$ws.Cells.Item(134, 8).Value = 2984.8572  # H134 was 3069.4814
$ws.Cells.Item(134, 9).Value = 2318.5833  # I134 was 2388.9565
$ws.Cells.Item(134, 11).Value = 6955.749899999999  # K134 was 7166.869499999999
$ws.Cells.Item(134, 13).Value = -4420.749899999999  # M134 was -4631.869499999999
$ws.Cells.Item(136, 8).Value = 1480.8182  # H136 was 2898.5454
$ws.Cells.Item(136, 9).Value = 1480.8182  # I136 was 1543.7778
$ws.Cells.Item(136, 10).Value = 0  # J136 was 8995
$ws.Cells.Item(136, 11).Value = 4442.4546  # K136 was 4631.3334
$ws.Cells.Item(136, 12).Value = 0  # L136 was 26985
$ws.Cells.Item(136, 13).Value = -1892.4546  # M136 was -2081.3334
$ws.Cells.Item(136, 14).ClearContents()  # N136 was -32085

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 1101.625  # H103 was 1232.2
$ws.Cells.Item(103, 9).Value = 0  # I103 was 599
$ws.Cells.Item(103, 10).Value = 1101.625  # J103 was 1302.5555
$ws.Cells.Item(103, 11).Value = 0  # K103 was 1797
$ws.Cells.Item(103, 12).Value = 3304.875  # L103 was 3907.6665
$ws.Cells.Item(103, 13).ClearContents()  # M103 was -918
$ws.Cells.Item(103, 14).Value = -5062.875  # N103 was -5665.666499999999
$ws.Cells.Item(113, 8).Value = 1672.9445  # H113 was 1687.5264
$ws.Cells.Item(113, 10).Value = 2085  # J113 was 2073.75
$ws.Cells.Item(113, 12).Value = 6255  # L113 was 6221.25
$ws.Cells.Item(113, 14).Value = -10595  # N113 was -10561.25
$ws.Cells.Item(136, 8).Value = 10322.333  # H136 was 10328.333
$ws.Cells.Item(136, 9).Value = 8983.5  # I136 was 8992.5
$ws.Cells.Item(136, 11).Value = 26950.5  # K136 was 26977.5
$ws.Cells.Item(136, 13).Value = -21850.5  # M136 was -21877.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 33750  # H15 was 40000
$ws.Cells.Item(15, 9).Value = 15000  # I15 was 0
$ws.Cells.Item(15, 11).Value = 15000  # K15 was 0
$ws.Cells.Item(15, 13).Value = -14712  # M15 was None
$ws.Cells.Item(81, 8).Value = 33750  # H81 was 40000
$ws.Cells.Item(81, 9).Value = 15000  # I81 was 0
$ws.Cells.Item(81, 11).Value = 15000  # K81 was 0
$ws.Cells.Item(81, 13).Value = -14002  # M81 was None
$ws.Cells.Item(84, 8).Value = 33750  # H84 was 40000
$ws.Cells.Item(84, 9).Value = 15000  # I84 was 0
$ws.Cells.Item(84, 11).Value = 45000  # K84 was 0
$ws.Cells.Item(84, 13).Value = -40008  # M84 was None
$ws.Cells.Item(122, 8).Value = 1750  # H122 was 3000
$ws.Cells.Item(122, 9).Value = 1750  # I122 was 3000
$ws.Cells.Item(122, 11).Value = 5250  # K122 was 9000
$ws.Cells.Item(122, 13).Value = -2800  # M122 was -6550

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(6, 8).Value = 0  # H6 was 20000
$ws.Cells.Item(6, 10).Value = 0  # J6 was 20000
$ws.Cells.Item(6, 12).Value = 0  # L6 was 20000
$ws.Cells.Item(6, 14).ClearContents()  # N6 was -20224
$ws.Cells.Item(22, 8).Value = 1685.7142  # H22 was 1344.2858
$ws.Cells.Item(22, 9).Value = 1100  # I22 was 477.5
$ws.Cells.Item(22, 10).Value = 2125  # J22 was 2500
$ws.Cells.Item(22, 11).Value = 1100  # K22 was 477.5
$ws.Cells.Item(22, 12).Value = 2125  # L22 was 2500
$ws.Cells.Item(22, 13).Value = -805  # M22 was -182.5
$ws.Cells.Item(22, 14).Value = -2715  # N22 was -3090
$ws.Cells.Item(27, 8).Value = 1685.7142  # H27 was 1344.2858
$ws.Cells.Item(27, 9).Value = 1100  # I27 was 477.5
$ws.Cells.Item(27, 10).Value = 2125  # J27 was 2500
$ws.Cells.Item(27, 11).Value = 1100  # K27 was 477.5
$ws.Cells.Item(27, 12).Value = 2125  # L27 was 2500
$ws.Cells.Item(27, 13).Value = -993  # M27 was -370.5
$ws.Cells.Item(27, 14).Value = -2339  # N27 was -2714
$ws.Cells.Item(96, 8).Value = 69393  # H96 was 69089.5
$ws.Cells.Item(122, 8).Value = 2835.182  # H122 was 2835.6365
$ws.Cells.Item(122, 9).Value = 2818.8  # I122 was 2977
$ws.Cells.Item(122, 10).Value = 2999  # J122 was 2199.5
$ws.Cells.Item(122, 11).Value = 8456.400000000001  # K122 was 8931
$ws.Cells.Item(122, 12).Value = 8997  # L122 was 6598.5
$ws.Cells.Item(122, 13).Value = -6006.400000000001  # M122 was -6481
$ws.Cells.Item(122, 14).Value = -13897  # N122 was -11498.5
$ws.Cells.Item(132, 8).Value = 5011.2085  # H132 was 5665.2856
$ws.Cells.Item(132, 9).Value = 5214.609  # I132 was 5931.9
$ws.Cells.Item(132, 11).Value = 15643.827  # K132 was 17795.7
$ws.Cells.Item(132, 13).Value = -13113.827  # M132 was -15265.7
$ws.Cells.Item(136, 8).Value = 1469.1666  # H136 was 1065.25
$ws.Cells.Item(136, 9).Value = 1469.1666  # I136 was 1065.25
$ws.Cells.Item(136, 13).Value = -1857.4998  # M136 was -645.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 974.1579  # H132 was 1005.94446
$ws.Cells.Item(132, 9).Value = 998  # I132 was 1033.0588
$ws.Cells.Item(132, 11).Value = 2994  # K132 was 3099.1764
$ws.Cells.Item(132, 13).Value = -464  # M132 was -569.1764000000003
$ws.Cells.Item(136, 8).Value = 1922.5834  # H136 was 2142.1667
$ws.Cells.Item(136, 9).Value = 1463.9584  # I136 was 1523.6957
$ws.Cells.Item(136, 10).Value = 2839.8333  # J136 was 3236.3845
$ws.Cells.Item(136, 11).Value = 4391.8752  # K136 was 4571.0871
$ws.Cells.Item(136, 12).Value = 8519.499899999999  # L136 was 9709.1535
$ws.Cells.Item(136, 13).Value = -1841.8752  # M136 was -2021.0871
$ws.Cells.Item(136, 14).Value = -13619.4999  # N136 was -14809.1535
